$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1034.3077
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H92").Value = 293.6111
$ws.Range("I92").Value = 185.14285
$ws.Range("K92").Value = 185.14285
$ws.Range("M92").Value = 1062.85715
$ws.Range("H107").Value = 1225.1305
$ws.Range("I107").Value = 1371.0555
$ws.Range("K107").Value = 1371.0555
$ws.Range("M107").Value = 548.9445000000001
$ws.Range("H111").Value = 1935.5
$ws.Range("I111").Value = 1080.6666
$ws.Range("J111").Value = 4500
$ws.Range("K111").Value = 3241.9998
$ws.Range("L111").Value = 13500
$ws.Range("M111").Value = -174.9998000000001
$ws.Range("N111").Value = -19634
$ws.Range("H137").Value = 1913.2
$ws.Range("I137").Value = 1894
$ws.Range("K137").Value = 5682
$ws.Range("M137").Value = -3132

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1324.909
$ws.Range("I74").Value = 933
$ws.Range("J74").Value = 2370
$ws.Range("K74").Value = 933
$ws.Range("L74").Value = 2370
$ws.Range("M74").Value = -59
$ws.Range("N74").Value = -4118
$ws.Range("H77").Value = 1324.909
$ws.Range("I77").Value = 933
$ws.Range("J77").Value = 2370
$ws.Range("K77").Value = 4665
$ws.Range("L77").Value = 11850
$ws.Range("M77").Value = -297
$ws.Range("N77").Value = -20586
$ws.Range("H110").Value = 1928.0646
$ws.Range("I110").Value = 1177.3684
$ws.Range("K110").Value = 1177.3684
$ws.Range("M110").Value = 867.6315999999999
$ws.Range("H132").Value = 4194.8
$ws.Range("I132").Value = 4052.1177
$ws.Range("K132").Value = 12156.3531
$ws.Range("M132").Value = -9626.3531

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 718.6667
$ws.Range("I5").Value = 718.6667
$ws.Range("K5").Value = 718.6667
$ws.Range("M5").Value = -605.6667
$ws.Range("H7").Value = 11111915
$ws.Range("I7").Value = 16666901
$ws.Range("J7").Value = 1942.6666
$ws.Range("K7").Value = 16666901
$ws.Range("L7").Value = 1942.6666
$ws.Range("M7").Value = -16666788
$ws.Range("N7").Value = -2168.6666
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22:N22").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4022.5833
$ws.Range("I31").Value = 2613.2856
$ws.Range("K31").Value = 2613.2856
$ws.Range("M31").Value = -2318.2856
$ws.Range("H34").Value = 4022.5833
$ws.Range("I34").Value = 2613.2856
$ws.Range("K34").Value = 2613.2856
$ws.Range("M34").Value = -2411.2856
$ws.Range("H132").Value = 1677.4445
$ws.Range("I132").Value = 1637.25
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 4911.75
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -2381.75
$ws.Range("N132").Value = -11057

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3697415.8
$ws.Range("I4").Value = 4692696
$ws.Range("J4").Value = 659.5714
$ws.Range("K4").Value = 14078088
$ws.Range("L4").Value = 1978.7142
$ws.Range("M4").Value = -14077976
$ws.Range("N4").Value = -2202.7142
$ws.Range("H14").Value = 745.9048
$ws.Range("I14").Value = 745.9048
$ws.Range("K14").Value = 2237.7144
$ws.Range("M14").Value = -2064.7144

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 65138.062
$ws.Range("I122").Value = 1990.091
$ws.Range("K122").Value = 5970.272999999999
$ws.Range("M122").Value = -3520.272999999999
$ws.Range("H132").Value = 2502.718
$ws.Range("I132").Value = 2672.9119
$ws.Range("K132").Value = 8018.7357
$ws.Range("M132").Value = -5488.7357

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3394.394
$ws.Range("I22").Value = 1947.875
$ws.Range("J22").Value = 4755.8237
$ws.Range("K22").Value = 1947.875
$ws.Range("L22").Value = 4755.8237
$ws.Range("M22").Value = -1652.875
$ws.Range("N22").Value = -5345.8237
$ws.Range("H27").Value = 3394.394
$ws.Range("I27").Value = 1947.875
$ws.Range("J27").Value = 4755.8237
$ws.Range("K27").Value = 1947.875
$ws.Range("L27").Value = 4755.8237
$ws.Range("M27").Value = -1840.875
$ws.Range("N27").Value = -4969.8237
$ws.Range("H55").Value = 533.6
$ws.Range("I55").Value = 502.46155
$ws.Range("J55").Value = 591.4286
$ws.Range("K55").Value = 502.46155
$ws.Range("L55").Value = 591.4286
$ws.Range("M55").Value = -329.46155
$ws.Range("N55").Value = -937.4286
$ws.Range("H93").Value = 537.25
$ws.Range("I93").Value = 374.5
$ws.Range("K93").Value = 374.5
$ws.Range("M93").Value = 873.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 14641.857
$ws.Range("I47").Value = 9000
$ws.Range("J47").Value = 15582.167
$ws.Range("K47").Value = 9000
$ws.Range("L47").Value = 15582.167
$ws.Range("M47").Value = -8428
$ws.Range("N47").Value = -16726.167
$ws.Range("H55").Value = 10000
$ws.Range("I55").Value = 10000
$ws.Range("J55").Value = 10000
$ws.Range("K55").Value = 10000
$ws.Range("L55").Value = 10000
$ws.Range("M55").Value = -9723
$ws.Range("N55").Value = -10554
$ws.Range("H58").Value = 60000
$ws.Range("I58").Value = 60000
$ws.Range("K58").Value = 60000
$ws.Range("M58").Value = -59692
$ws.Range("H64").Value = 63330
$ws.Range("J64").Value = 63500
$ws.Range("L64").Value = 63500
$ws.Range("N64").Value = -63996
$ws.Range("H67").Value = 63330
$ws.Range("J67").Value = 63500
$ws.Range("L67").Value = 63500
$ws.Range("N67").Value = -65216
$ws.Range("H100").Value = 1577.25
$ws.Range("I100").Value = 1683
$ws.Range("K100").Value = 3366
$ws.Range("M100").Value = -2825
$ws.Range("H129").Value = 81899
$ws.Range("J129").Value = 81899
$ws.Range("L129").Value = 81899
$ws.Range("N129").Value = -91899
$ws.Range("H132").Value = 2632.7778
$ws.Range("I132").Value = 2336.875
$ws.Range("K132").Value = 7010.625
$ws.Range("M132").Value = -4480.625
